# Update the ER-diagram on the "Data Model" slide (slide 11): the PRODUCT and
# LINE_ITEM entity boxes (header + field list), their relationship-label boxes,
# and the two bent connectors swap places/labels/geometry.

$eps = 0.00002   # nudge so Single-precision EMU round-trip lands on the exact value

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(11)

# --- Entity headers -------------------------------------------------------
# "LINE_ITEM" header box now reads "PRODUCT"
$s.Shapes.Item("Shape 12").TextFrame.TextRange.Text = "PRODUCT"

# "LINE_ITEM" field list becomes the PRODUCT field list
$s.Shapes.Item("Shape 13").TextFrame.TextRange.Text = "int id PK`nstring name`ndecimal price`nint stock"

# --- Relationship labels ----------------------------------------------------
# "ntains" label becomes the order-side label
$s.Shapes.Item("Shape 16").TextFrame.TextRange.Text = 'rdered in"'

# "PRODUCT" header box now reads "ntains"
$s.Shapes.Item("Shape 18").TextFrame.TextRange.Text = "ntains"

# PRODUCT field box shrinks to an empty label box
$shape19 = $s.Shapes.Item("Shape 19")
$shape19.Height = (280000 / 12700) + $eps
$shape19.TextFrame.TextRange.Text = ""

# the other relationship label box now reads "LINE_ITEM"
$s.Shapes.Item("Shape 20").TextFrame.TextRange.Text = "LINE_ITEM"

# empty field box grows and gets the LINE_ITEM field list
$shape21 = $s.Shapes.Item("Shape 21")
$shape21.Height = (1120000 / 12700) + $eps
$shape21.TextFrame.TextRange.Text = "int id PK`nint order_id FK`nint product_id FK`nint quantity"

# --- Connectors --------------------------------------------------------------
# Connector 62: was flipped, short diagonal under ORDER -> becomes a short,
# unflipped vertical run between the two left-hand boxes.
$c62 = $s.Shapes.Item("Connector 62")
$c62.HorizontalFlip = 0
$c62.Left = (2700000 / 12700) + $eps
$c62.Top = (1800000 / 12700) + $eps
$c62.Width = (600000 / 12700) + $eps
$c62.Height = (2500000 / 12700) + $eps

# Connector 63: was a short flat stub -> becomes the long flipped diagonal
# spanning from ORDER down to PRODUCT.
$c63 = $s.Shapes.Item("Connector 63")
$c63.HorizontalFlip = -1
$c63.Left = (500000 / 12700) + $eps
$c63.Top = (1800000 / 12700) + $eps
$c63.Width = (5000000 / 12700) + $eps
$c63.Height = (2500000 / 12700) + $eps
